# Applies the diff:
#  - Column C ("Förändrad" / Modified date) for data rows 2..407:
#      value 45184 -> 45186
#  - HYPERLINK formulas in columns S, T, V, W, X, Y for data rows 2..14:
#      add a second argument to HYPERLINK() equal to the row's "Beteckning"
#      (column A) value, e.g.
#        HYPERLINK("...url...")  ->  HYPERLINK("...url...", "A 31390-2020")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 407
$lastHyperlinkRow = 14

# Column indexes: A=1, C=3, S=19, T=20, V=22, W=23, X=24, Y=25
$colBeteckning = 1
$colForandrad = 3
$hyperlinkCols = 19, 20, 22, 23, 24, 25

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    # Update the "Forandrad" date value from 45184 to 45186
    $cModified = $ws.Cells.Item($row, $colForandrad)
    if ($cModified.Value2 -eq 45184) {
        $cModified.Value2 = 45186
    }

    # Only the first 13 data rows (2..14) have hyperlink formulas
    if ($row -le $lastHyperlinkRow) {
        $beteckning = $ws.Cells.Item($row, $colBeteckning).Value2

        foreach ($col in $hyperlinkCols) {
            $cell = $ws.Cells.Item($row, $col)
            $formula = $cell.Formula
            if ($formula -and $formula.Length -gt 0) {
                $closingParenIndex = $formula.LastIndexOf(")")
                $alreadyHasSecondArg = ($formula.IndexOf(", """) -ge 0) -or ($formula.IndexOf(",""") -ge 0)
                if ($closingParenIndex -ge 0 -and -not $alreadyHasSecondArg) {
                    $newFormula = $formula.Substring(0, $closingParenIndex) + ', "' + $beteckning + '")'
                    $cell.Formula = $newFormula
                }
            }
        }
    }
}

Write-Host "done"
